$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Date_1 (column A) needs updating to 2025/11/29
$rows = @(2, 8, 14, 20, 26, 32, 38, 44, 50, 56, 62, 68, 74)

# New EBITDA (column B) values for the rows that also changed
$bvals = @{
    2  = "5.11"
    8  = "7.67"
    14 = "2.83"
    20 = "12.37"
    26 = "9.97"
    32 = "26.09"
    44 = "10.97"
    56 = "35.94"
    62 = "11.45"
    68 = "12.37"
    74 = "16.06"
}

foreach ($r in $rows) {
    # Use a pristine same-row cell (column C, untouched, default style) as a
    # style donor so the text stays plain text (no stray date/number format).
    $styleDonor = $ws.Cells.Item($r, 3)

    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = "2025/11/29"
    $cellA.Style = $styleDonor.Style

    if ($bvals.ContainsKey($r)) {
        $cellB = $ws.Cells.Item($r, 2)
        $cellB.NumberFormat = "@"
        $cellB.Value = $bvals[$r]
        $cellB.Style = $styleDonor.Style
    }
}
